$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear cells that no longer hold values
$ws.Range("E2").ClearContents()
$ws.Range("E3").ClearContents()
$ws.Range("F4").ClearContents()
$ws.Range("G5").ClearContents()
$ws.Range("G6").ClearContents()

# Update row 8 values
$ws.Range("E8").Value = "https://quizizz.com/join?gc=08539312"
$ws.Range("F8").Value = "/forest.pdf"
$ws.Range("G8").Value = "fdfdfd"

# Update selection to match the final view state
$ws.Range("E2:G6").Select()
